$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 8.898150666666666
$ws.Range("H2").Value = 26.694452
$ws.Range("I2").Value = 0.3765197173862137
$ws.Range("J2").Value = 0.3765197173862137
$ws.Range("M2").Value = 0.008424666666666667
$ws.Range("N2").Value = 0.025274
$ws.Range("O2").Value = 0.0215396310803394
$ws.Range("P2").Value = 0.0215396310803394
$ws.Range("Q2").Value = 0.07496395331644444
$ws.Range("R2").Value = 0.674675579848
$ws.Range("S2").Value = 0.008110095806972696
$ws.Range("T2").Value = 0.008110095806972696
$ws.Range("G3").Value = 8.898150666666666
$ws.Range("H3").Value = 26.694452
$ws.Range("I3").Value = 0.3765197173862137
$ws.Range("J3").Value = 0.3765197173862137
$ws.Range("O3").Value = 0.9647613885451503
$ws.Range("P3").Value = 0.9647613885451503
$ws.Range("Q3").Value = 3.357640036760889
$ws.Range("R3").Value = 30.218760330848
$ws.Range("S3").Value = 0.3632516853601511
$ws.Range("T3").Value = 0.3632516853601511
$ws.Range("G4").Value = 8.898150666666666
$ws.Range("H4").Value = 26.694452
$ws.Range("I4").Value = 0.3765197173862137
$ws.Range("J4").Value = 0.3765197173862137
$ws.Range("M4").Value = 0.005357999999999999
$ws.Range("N4").Value = 0.016074
$ws.Range("O4").Value = 0.01369898037451038
$ws.Range("P4").Value = 0.01369898037451039
$ws.Range("Q4").Value = 0.04767629127199999
$ws.Range("R4").Value = 0.4290866214479999
$ws.Range("S4").Value = 0.005157936219089938
$ws.Range("T4").Value = 0.005157936219089938
$ws.Range("I5").Value = 0.1415167724465014
$ws.Range("J5").Value = 0.1415167724465015
$ws.Range("M5").Value = 0.008424666666666667
$ws.Range("N5").Value = 0.025274
$ws.Range("O5").Value = 0.0215396310803394
$ws.Range("P5").Value = 0.0215396310803394
$ws.Range("Q5").Value = 0.02817556752888889
$ws.Range("R5").Value = 0.2535801077600001
$ws.Range("S5").Value = 0.003048219070177981
$ws.Range("T5").Value = 0.003048219070177981
$ws.Range("I6").Value = 0.1415167724465014
$ws.Range("J6").Value = 0.1415167724465015
$ws.Range("O6").Value = 0.9647613885451503
$ws.Range("P6").Value = 0.9647613885451503
$ws.Range("S6").Value = 0.1365299178879148
$ws.Range("T6").Value = 0.1365299178879148
$ws.Range("I7").Value = 0.1415167724465014
$ws.Range("J7").Value = 0.1415167724465015
$ws.Range("M7").Value = 0.005357999999999999
$ws.Range("N7").Value = 0.016074
$ws.Range("O7").Value = 0.01369898037451038
$ws.Range("P7").Value = 0.01369898037451039
$ws.Range("Q7").Value = 0.01791936664
$ws.Range("R7").Value = 0.16127429976
$ws.Range("S7").Value = 0.001938635488408675
$ws.Range("T7").Value = 0.001938635488408676
$ws.Range("G8").Value = 11.39006466666667
$ws.Range("H8").Value = 34.170194
$ws.Range("I8").Value = 0.4819635101672848
$ws.Range("J8").Value = 0.4819635101672848
$ws.Range("M8").Value = 0.008424666666666667
$ws.Range("N8").Value = 0.025274
$ws.Range("O8").Value = 0.0215396310803394
$ws.Range("P8").Value = 0.0215396310803394
$ws.Range("Q8").Value = 0.09595749812844445
$ws.Range("R8").Value = 0.8636174831560001
$ws.Range("S8").Value = 0.01038131620318872
$ws.Range("T8").Value = 0.01038131620318872
$ws.Range("G9").Value = 11.39006466666667
$ws.Range("H9").Value = 34.170194
$ws.Range("I9").Value = 0.4819635101672848
$ws.Range("J9").Value = 0.4819635101672848
$ws.Range("O9").Value = 0.9647613885451503
$ws.Range("P9").Value = 0.9647613885451503
$ws.Range("Q9").Value = 4.297942188072889
$ws.Range("R9").Value = 38.681479692656
$ws.Range("S9").Value = 0.4649797852970843
$ws.Range("T9").Value = 0.4649797852970843
$ws.Range("G10").Value = 11.39006466666667
$ws.Range("H10").Value = 34.170194
$ws.Range("I10").Value = 0.4819635101672848
$ws.Range("J10").Value = 0.4819635101672848
$ws.Range("M10").Value = 0.005357999999999999
$ws.Range("N10").Value = 0.016074
$ws.Range("O10").Value = 0.01369898037451038
$ws.Range("P10").Value = 0.01369898037451039
$ws.Range("Q10").Value = 0.061027966484
$ws.Range("R10").Value = 0.549251698356
$ws.Range("S10").Value = 0.00660240866701177
$ws.Range("T10").Value = 0.006602408667011772
